$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Archetypes": insert a new column E ("source": w/s) between the
# existing "size code" column (D) and the "description" column (old E, which
# becomes F); retitle D1 "code" -> "size"; drop the last data row (row 11),
# keeping only the empty, styled B11 cell.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Archetypes")

# Insert a new column before the current column E (pushes description -> F,
# and the lone I5 style cell -> J5).
$ws.Columns.Item(5).Insert()
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(4).ColumnWidth

# Header row
$ws.Range("D1").Value = "size"
$ws.Range("E1").Value = "source"
$ws.Range("F1").Value = "description"

# New "source" column values (w = wind, s = solar)
$ws.Range("E2").Value = "w"
$ws.Range("E3").Value = "s"
$ws.Range("E4").Value = "w"
$ws.Range("E5").Value = "w"
$ws.Range("E6").Value = "s"
$ws.Range("E7").Value = "w"
$ws.Range("E8").Value = "w"
$ws.Range("E9").Value = "s"
$ws.Range("E10").Value = "w"

# Drop the former "solar generating district" row, leaving only the styled,
# empty B11 cell behind.
$ws.Range("A11").Value = ""
$ws.Range("B11").Value = ""
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = ""

$ws.Range("A2:A10").Select()

# ---------------------------------------------------------------------------
# Navigation / view-state tweaks captured in the diff (harmless no-ops on
# content, but keep the workbook's saved UI state in sync with the source).
# ---------------------------------------------------------------------------
$wsOffshore = $wb.Worksheets.Item("Offshore wind")
$wsOffshore.Range("I24").Select()

$wsOnshore = $wb.Worksheets.Item("Onshore wind")
$wsOnshore.Range("D31").Select()

$wsSolar = $wb.Worksheets.Item("Solar Photovoltaic")
$wsSolar.Range("H37").Select()

$ws.Select()
